# The commit removes the "F-statistic" row's computed Treatment/Control
# values (C42 = "3.86", D42 = "0.356") from the balance table, which were
# the joint-test F-statistic and its p-value. Clearing these two cells
# also drops their now-unused entries from the shared string table
# (uniqueCount 146 -> 144), and the engine automatically re-packs/reindexes
# the remaining shared strings, which is what shifts every other shared
# string reference in column D throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C42").Value = ""
$ws.Range("D42").Value = ""
